$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 2749  # was 2748
$wsExhibit.Range("F5").Value = 943  # was 942
$wsExhibit.Range("F7").Value = 2401  # was 2395
$wsExhibit.Range("F8").Value = 1858  # was 1856
$wsExhibit.Range("F11").Value = 2511  # was 2510
$wsExhibit.Range("F13").Value = 253  # was 252
$wsExhibit.Range("F17").Value = 122  # was 121
$wsExhibit.Range("F18").Value = 9372  # was 9364
$wsExhibit.Range("F19").Value = 61  # was 59
$wsExhibit.Range("F20").Value = 7294  # was 7287
$wsExhibit.Range("F21").Value = 11854  # was 11844
$wsExhibit.Range("F25").Value = 370  # was 369
$wsExhibit.Range("F26").Value = 568  # was 567
$wsExhibit.Range("F27").Value = 2657  # was 2653
$wsExhibit.Range("F29").Value = 206  # was 205
$wsExhibit.Range("F30").Value = 2614  # was 2612
$wsExhibit.Range("F31").Value = 814  # was 804
$wsExhibit.Range("F32").Value = 53  # was 51
$wsExhibit.Range("F33").Value = 4538  # was 4537
$wsExhibit.Range("F34").Value = 994  # was 990
$wsExhibit.Range("F37").Value = 548  # was 547

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F6").Value = 16  # was 15
$wsShow.Range("F19").Value = 1  # was 0

$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F4").Value = 167  # was 166

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 2749  # was 2748
$wsAll.Range("F8").Value = 943  # was 942
$wsAll.Range("F11").Value = 2401  # was 2395
$wsAll.Range("F12").Value = 16  # was 15
$wsAll.Range("F13").Value = 1858  # was 1856
$wsAll.Range("F15").Value = 2511  # was 2510
$wsAll.Range("F18").Value = 253  # was 252
$wsAll.Range("F21").Value = 122  # was 121
$wsAll.Range("F22").Value = 9372  # was 9364
$wsAll.Range("F23").Value = 61  # was 59
$wsAll.Range("F24").Value = 7294  # was 7287
$wsAll.Range("F25").Value = 11854  # was 11844
$wsAll.Range("F29").Value = 370  # was 369
$wsAll.Range("F31").Value = 568  # was 567
$wsAll.Range("F33").Value = 2657  # was 2653
$wsAll.Range("F37").Value = 206  # was 205
$wsAll.Range("F38").Value = 53  # was 51
$wsAll.Range("F39").Value = 4538  # was 4537
$wsAll.Range("F41").Value = 1  # was 0
$wsAll.Range("F46").Value = 548  # was 547
